# Insert a new data row before current row 102 (shifts existing rows 102-211 down
# to 103-212) and populate the new row with the latest weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(102).Insert()

$ws.Cells.Item(102, 1).Value = 7
$ws.Cells.Item(102, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(102, 3).Value = "Ñuble"
$ws.Cells.Item(102, 4).Value = 44895
$ws.Cells.Item(102, 5).Value = 16
$ws.Cells.Item(102, 6).Value = 100112045
$ws.Cells.Item(102, 7).Value = "Zapallo"
$ws.Cells.Item(102, 8).Value = "Camote"
$ws.Cells.Item(102, 9).Value = "1a nueva(o)"
$ws.Cells.Item(102, 10).Value = 400
$ws.Cells.Item(102, 11).Value = 1100
$ws.Cells.Item(102, 12).Value = 1200
$ws.Cells.Item(102, 13).Value = 1150
$ws.Cells.Item(102, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(102, 15).Value = "Perú"
$ws.Cells.Item(102, 16).Value = 1150
$ws.Cells.Item(102, 17).Value = 1
$ws.Cells.Item(102, 18).Value = "Hortaliza"
